$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.93
$ws.Range("H2").Value = 4.3
$ws.Range("G3").Value = 2.66
$ws.Range("I3").Value = 3.75
$ws.Range("J3").Value = 3.05
$ws.Range("Q3").Value = 2.2
$ws.Range("S3").Value = 3.85
$ws.Range("V3").Value = 1.34
$ws.Range("W3").Value = 1.61
$ws.Range("J5").Value = 3.95
$ws.Range("AO5").Value = 55
$ws.Range("H6").Value = 4.6
$ws.Range("I6").Value = 4.7
$ws.Range("P6").Value = 1.57
$ws.Range("V6").Value = 1.27
$ws.Range("G7").Value = 2.26
$ws.Range("W7").Value = 1.79
$ws.Range("G8").Value = 2.54
$ws.Range("H8").Value = 2.72
$ws.Range("K8").Value = 4.7
$ws.Range("S8").Value = 1.95
$ws.Range("W8").Value = 1.65
$ws.Range("G9").Value = 2.26
$ws.Range("I9").Value = 3.3
$ws.Range("N9").Value = 6.2
$ws.Range("Q9").Value = 1.36
$ws.Range("V9").Value = 1.43
$ws.Range("W9").Value = 1.79
$ws.Range("G10").Value = 1.95
$ws.Range("H10").Value = 3.8
$ws.Range("I10").Value = 4.3
$ws.Range("L10").Value = 1.2
$ws.Range("P10").Value = 2.9
$ws.Range("V10").Value = 1.3
$ws.Range("G11").Value = 2.28
$ws.Range("J11").Value = 3.15
$ws.Range("K11").Value = 3.7
$ws.Range("L11").Value = 1.37
$ws.Range("O11").Value = 1.38
$ws.Range("S11").Value = 3.75
$ws.Range("W11").Value = 1.78
$ws.Range("X11").Value = 14.5
$ws.Range("AM11").Value = 140
$ws.Range("I12").Value = 2.84
$ws.Range("N12").Value = 2.9
$ws.Range("Q12").Value = 2.38
$ws.Range("R12").Value = 1.23
$ws.Range("S12").Value = 4.6
$ws.Range("Y12").Value = 9
$ws.Range("F13").Value = 1.26
$ws.Range("H13").Value = 15
$ws.Range("Q13").Value = 1.69
$ws.Range("V13").Value = 1.06
$ws.Range("AD13").Value = 60
$ws.Range("F14").Value = 2.86
$ws.Range("G14").Value = 2.9
$ws.Range("P14").Value = 1.82
$ws.Range("W14").Value = 1.52
$ws.Range("X14").Value = 12
